$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.101.43"
$ws.Range("E2").Value = "  +5.28%  "
$ws.Range("D3").Value = "'1.920.33"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  -1.05%  "
$ws.Range("D5").Value = "'326.83"
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").Value = "'0.5164"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").Value = "'0.4013"
$ws.Range("E8").Value = "  +2.76%  "
$ws.Range("D9").Value = "'0.08451"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.122"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'42.62"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "'21.68"
$ws.Range("E12").Value = "  +6.02%  "
$ws.Range("D13").Value = "'6.340"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "'1.923.84"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "'7.349"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "'96.04"
$ws.Range("E17").Value = "  +5.11%  "
$ws.Range("D18").Value = "'0.00001115"
$ws.Range("D19").Value = "'0.06727"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "'6.057"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").Value = "'30.101.92"
$ws.Range("E23").Value = "  +5.11%  "
$ws.Range("D24").Value = "'11.23"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").Value = "'2.202"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").Value = "'2.142.07"
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("D27").Value = "'160.72"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "'21.01"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "'2.456"
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").Value = "'128.86"
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "'1.073"
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("D32").Value = "'0.1059"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").Value = "'6.079"
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("D34").Value = "'3.660"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").Value = "'0.02514"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("D36").Value = "'0.06595"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "'0.2217"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("D38").Value = "'1.238"
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").Value = "'9.007"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("D40").Value = "'5.200"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").Value = "'0.6545"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "'1.244"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").Value = "'11.42"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("D44").Value = "'0.6133"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").Value = "'13.13"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'3.756"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").Value = "'2.053"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").Value = "'1.241"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'125.45"
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("D51").Value = "'79.26"
$ws.Range("E51").Value = "  +3.19%  "
